$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.337909340858459
$ws.Range("B1").Value = 3.299721002578735
$ws.Range("C1").Value = 5.521449089050293
$ws.Range("D1").Value = 1.695821166038513
$ws.Range("E1").Value = 0.9910944700241089
